## Add files via upload
## - Adds a new "Presets" worksheet at the end of the workbook with a small
##   lookup table of preset column combinations.
## - Updates the window/selection state: the "Mandatory" sheet keeps a
##   selection of V2 (and loses its special top-left/selected-tab state),
##   while the new "Presets" sheet becomes the active/selected tab with
##   B6 selected.

$wb = $excel.ActiveWorkbook

# --- capture source sheet used to copy the existing header/gray style ---
$styleSource = $wb.Worksheets.Item("Input").Range("A1")

# --- make sure "Mandatory" ends up with a plain selection of V2 (no more
#     topLeftCell scroll position, no more tab-selected flag) ---
$wsMandatory = $wb.Worksheets.Item("Mandatory")
$wsMandatory.Activate() | Out-Null
$wsMandatory.Range("V2").Select() | Out-Null

# --- add the new "Presets" sheet after the last existing sheet ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsPresets = $wb.Worksheets.Add($null, $lastSheet)
$wsPresets.Name = "Presets"

# --- header row (written in MSD1, MSD2, Mcomb, IPD, EsP, EsC, EsD order so
#     new shared-string entries land in that sequence, matching the source
#     workbook's string table) ---
$wsPresets.Cells.Item(1, 1).Value = "MSD1"
$wsPresets.Cells.Item(1, 2).Value = "MSD2"
$wsPresets.Cells.Item(1, 6).Value = "Mcomb"
$wsPresets.Cells.Item(1, 7).Value = "IPD"
$wsPresets.Cells.Item(1, 3).Value = "EsP"
$wsPresets.Cells.Item(1, 4).Value = "EsC"
$wsPresets.Cells.Item(1, 5).Value = "EsD"

# --- numeric rows 2-4 ---
$wsPresets.Cells.Item(2, 1).Value = 1
$wsPresets.Cells.Item(2, 2).Value = 1
$wsPresets.Cells.Item(2, 3).Value = 2
$wsPresets.Cells.Item(2, 4).Value = 2
$wsPresets.Cells.Item(2, 5).Value = 2
$wsPresets.Cells.Item(2, 6).Value = 3
$wsPresets.Cells.Item(2, 7).Value = 4

$wsPresets.Cells.Item(3, 1).Value = 1
$wsPresets.Cells.Item(3, 2).Value = 2
$wsPresets.Cells.Item(3, 3).Value = 1
$wsPresets.Cells.Item(3, 4).Value = 2
$wsPresets.Cells.Item(3, 5).Value = 2
$wsPresets.Cells.Item(3, 6).Value = 2
$wsPresets.Cells.Item(3, 7).Value = 1

$wsPresets.Cells.Item(4, 1).Value = 0
$wsPresets.Cells.Item(4, 2).Value = 0
$wsPresets.Cells.Item(4, 3).Value = 0
$wsPresets.Cells.Item(4, 4).Value = 0
$wsPresets.Cells.Item(4, 5).Value = 0
$wsPresets.Cells.Item(4, 6).Value = 0
$wsPresets.Cells.Item(4, 7).Value = 0

# --- row 5 : field-name labels (Study_ID across, patient_data in G) ---
$wsPresets.Cells.Item(5, 1).Value = "Study_ID"
$wsPresets.Cells.Item(5, 2).Value = "Study_ID"
$wsPresets.Cells.Item(5, 3).Value = "Study_ID"
$wsPresets.Cells.Item(5, 4).Value = "Study_ID"
$wsPresets.Cells.Item(5, 5).Value = "Study_ID"
$wsPresets.Cells.Item(5, 6).Value = "Study_ID"
$wsPresets.Cells.Item(5, 7).Value = "patient_data"

# --- row 6 ---
$wsPresets.Cells.Item(6, 1).Value = "Mean"
$wsPresets.Cells.Item(6, 3).Value = "N"
$wsPresets.Cells.Item(6, 4).Value = "group_ID"
$wsPresets.Cells.Item(6, 5).Value = "group_ID"
$wsPresets.Cells.Item(6, 6).Value = "Mean"
$wsPresets.Cells.Item(6, 7).Value = "labs"

# --- row 7 ---
$wsPresets.Cells.Item(7, 1).Value = "SD"
$wsPresets.Cells.Item(7, 3).Value = "N_total"
$wsPresets.Cells.Item(7, 4).Value = "Mean"
$wsPresets.Cells.Item(7, 5).Value = "N"
$wsPresets.Cells.Item(7, 6).Value = "SD"

# --- row 8 ---
$wsPresets.Cells.Item(8, 1).Value = "N"
$wsPresets.Cells.Item(8, 4).Value = "SD"
$wsPresets.Cells.Item(8, 5).Value = "N_total"
$wsPresets.Cells.Item(8, 6).Value = "N"

# --- row 9 ---
$wsPresets.Cells.Item(9, 1).Value = "SE"
$wsPresets.Cells.Item(9, 4).Value = "N"
$wsPresets.Cells.Item(9, 6).Value = "labs"

# --- rows 10-16, column A only ---
$wsPresets.Cells.Item(10, 1).Value = "Median"
$wsPresets.Cells.Item(11, 1).Value = "q1"
$wsPresets.Cells.Item(12, 1).Value = "q3"
$wsPresets.Cells.Item(13, 1).Value = "min"
$wsPresets.Cells.Item(14, 1).Value = "max"
$wsPresets.Cells.Item(15, 1).Value = "ulci"
$wsPresets.Cells.Item(16, 1).Value = "llci"

# --- apply the existing gray header style to the handful of highlighted
#     cells (copy-format only, so the underlying values are untouched) ---
$styleSource.Copy()
$wsPresets.Range("G5").PasteSpecial(-4122)
$wsPresets.Range("D6").PasteSpecial(-4122)
$wsPresets.Range("E6").PasteSpecial(-4122)
$wsPresets.Range("G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- final view state: Presets is the active sheet, with B6 selected ---
$wsPresets.Activate() | Out-Null
$wsPresets.Range("B6").Select() | Out-Null
